$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, date range) ---
$ws.Range("A8").Value = "Volume 29   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# --- Plain value updates (style unchanged) ---
$ws.Range("M15").Value = 15.789473684210
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 244
$ws.Range("J16").Value = 180
$ws.Range("K16").Value = 35.555555555555
$ws.Range("L16").Value = 53.459119496855
$ws.Range("M16").Value = 0.411522633744
$ws.Range("N16").Value = -79.717373233582
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 22.727272727272
$ws.Range("I17").Value = 272
$ws.Range("J17").Value = 246
$ws.Range("K17").Value = 10.569105691056
$ws.Range("L17").Value = 40.932642487046
$ws.Range("M17").Value = 39.487179487179
$ws.Range("N17").Value = -65.350318471337
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 8.333333333333
$ws.Range("I18").Value = 130
$ws.Range("J18").Value = 96
$ws.Range("K18").Value = 35.416666666666
$ws.Range("L18").Value = -32.989690721649
$ws.Range("M18").Value = -18.75
$ws.Range("N18").Value = -92.605233219567
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 77.777777777777
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 22.727272727272
$ws.Range("I19").Value = 505
$ws.Range("J19").Value = 504
$ws.Range("K19").Value = 0.198412698412
$ws.Range("L19").Value = 28.172588832487
$ws.Range("M19").Value = 56.346749226006
$ws.Range("N19").Value = -52.268431001890
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 263
$ws.Range("J20").Value = 178
$ws.Range("K20").Value = 47.752808988764
$ws.Range("L20").Value = 91.970802919708
$ws.Range("M20").Value = 237.179487179487
$ws.Range("N20").Value = -83.982947624847
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 6.666666666666
$ws.Range("F21").Value = 130
$ws.Range("H21").Value = 14.035087719298
$ws.Range("I21").Value = 1440
$ws.Range("J21").Value = 1228
$ws.Range("K21").Value = 17.263843648208
$ws.Range("L21").Value = 31.267092069279
$ws.Range("M21").Value = 40.900195694716
$ws.Range("N21").Value = -78.092195344591
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 66.666666666666
$ws.Range("I22").Value = 26
$ws.Range("J22").Value = 23
$ws.Range("K22").Value = 13.043478260869
$ws.Range("L22").Value = 73.333333333333
$ws.Range("M22").Value = 23.809523809523
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 35
$ws.Range("J23").Value = 28
$ws.Range("K23").Value = 25
$ws.Range("L23").Value = 16.666666666666
$ws.Range("M23").Value = 75
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 13.636363636363
$ws.Range("F24").Value = 101
$ws.Range("G24").Value = 88
$ws.Range("H24").Value = 14.772727272727
$ws.Range("I24").Value = 1219
$ws.Range("J24").Value = 767
$ws.Range("K24").Value = 58.930899608865
$ws.Range("L24").Value = 43.580683156654
$ws.Range("M24").Value = 119.244604316547
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -10
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 13.793103448275
$ws.Range("I25").Value = 386
$ws.Range("J25").Value = 357
$ws.Range("K25").Value = 8.123249299719
$ws.Range("L25").Value = 11.239193083573
$ws.Range("M25").Value = -18.736842105263
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 46
$ws.Range("J27").Value = 53
$ws.Range("K27").Value = -13.207547169811
$ws.Range("L27").Value = 17.948717948717
$ws.Range("N28").Value = -89.221556886227
$ws.Range("N29").Value = -87.671232876712
$ws.Range("I30").Value = 5
$ws.Range("K30").Value = 400
$ws.Range("L30").Value = 25

# --- Special updates requiring a style change (copy format from a donor cell) ---
$ws.Range("F15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("D17").Value = 7
$ws.Range("I14").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = -14.285714285714
$ws.Range("K14").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("C22").Value = "'0"
$ws.Range("F14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = 2
$ws.Range("J14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("L14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("D26").Value = "'0"
$ws.Range("G14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "'***.*"
$ws.Range("H14").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$excel.CutCopyMode = $false
